$wb = $excel.ActiveWorkbook

$wsCasting = $wb.Worksheets.Item("Casting Method")
$wsCasting.Select()
$wsCasting.Range("A8").Select()

$wsOnHit = $wb.Worksheets.Item("OnHit Effect")
$wsOnHit.Select()

$wsOnHit.Range("A12").Value = "Sleep"
$wsOnHit.Range("B12").Value = "Puts target to sleep state. Taking damage will aggro"

$wsOnHit.Range("A13").Value = "Bullrush"
$wsOnHit.Range("B13").Value = "Target will burst forward certain distance depending on power"

$wsOnHit.Range("A14").Value = "Boost Attack Power"
$wsOnHit.Range("A15").Value = "Boost Defense"
$wsOnHit.Range("A16").Value = "Spawn Minion"

$wsOnHit.Columns.Item(2).ColumnWidth = 56.6

$wsOnHit.Range("A14").Select()
